# Build site at 2022-09-26 16:07:08 UTC
# Reproduces the authored edit to LOQ4256.xlsx:
#  - The row that only held "5701460 - Antonio Iacono" in B/C (old row 13,
#    with no label in column A) is removed entirely, shifting every row
#    below it up by one.
#  - A handful of the label rows below keep their own label (column A) but
#    end up showing the content value that used to belong to a different
#    row, because the content columns (B/C) were not re-synced after the
#    row delete.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$docente = "5701460 - Antonio Iacono"
$dataAtivacao = "01/01/2021"
$metodoTexto = "Aulas Expositivas; trabalhos em grupo; exercícios individuais; palestras e painel integrado."
$criterioTexto = "MF = (0,40*P1 + 0,40*P2 + 0,20*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários."
$normaTexto = "NF = (MF + PR)/2, onde PR é uma prova de recuperação."

# 1) Remove the stray row (old row 13) that held only the docente
#    responsável value with no label -- everything below shifts up by one.
$ws.Rows.Item(13).Delete()

# 2) Row 10 ("Objetivos:") now shows the docente responsável text.
$ws.Range("B10").Value = $docente
$ws.Range("C10").Value = $docente

# 3) Row 13 ("Programa resumido:", was row 14) now reads "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# 4) Row 15 ("Programa:", was row 16) now shows the activation date.
#    Force this in as text (matching row 8's existing "01/01/2021" text
#    cell) instead of letting Excel auto-convert it to a date serial, by
#    pre-formatting as Text and then restoring the normal cell format
#    (copied from B8/C8, which already carry the same string + style).
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = $dataAtivacao
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = $dataAtivacao

$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4122)

# 5) Row 18 ("Método:", was row 19) now shows the docente responsável text.
$ws.Range("B18").Value = $docente
$ws.Range("C18").Value = $docente

# 6) Row 19 ("Critério:", was row 20) now shows the old Método text.
$ws.Range("B19").Value = $metodoTexto
$ws.Range("C19").Value = $metodoTexto

# 7) Row 20 ("Norma de recuperação:", was row 21) now shows the old
#    Critério text.
$ws.Range("B20").Value = $criterioTexto
$ws.Range("C20").Value = $criterioTexto

# 8) Row 21 ("Bibliografia:", was row 22) now shows the old Norma de
#    recuperação text.
$ws.Range("B21").Value = $normaTexto
$ws.Range("C21").Value = $normaTexto
